$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H54").Value = 299.25
$ws.Range("I54").Value = 299.25
$ws.Range("K54").Value = 299.25
$ws.Range("M54").Value = 186.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 564559.8
$ws.Range("J112").Value = 745736.1
$ws.Range("L112").Value = 2237208.3
$ws.Range("N112").Value = -2239424.3

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 27796986
$ws.Range("I137").Value = 8334282
$ws.Range("J137").Value = 52125370
$ws.Range("K137").Value = 25002846
$ws.Range("L137").Value = 156376110
$ws.Range("M137").Value = -25000296
$ws.Range("N137").Value = -156381210

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3644.8413
$ws.Range("I138").Value = 5261.9165
$ws.Range("J138").Value = 3264.353
$ws.Range("K138").Value = 15785.7495
$ws.Range("L138").Value = 9793.059000000001
$ws.Range("M138").Value = -10645.7495
$ws.Range("N138").Value = -20073.059

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9753.062
$ws.Range("I32").Value = 7701.6895
$ws.Range("K32").Value = 7701.6895
$ws.Range("M32").Value = -7414.6895

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H47").Value = 10000
$ws.Range("J47").Value = 10000
$ws.Range("L47").Value = 10000
$ws.Range("N47").Value = -11450

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 3111.111
$ws.Range("I102").Value = 2400
$ws.Range("K102").Value = 2400
$ws.Range("M102").Value = -778

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 21545840
$ws.Range("I132").Value = 24330614
$ws.Range("J132").Value = 8550230
$ws.Range("K132").Value = 72991842
$ws.Range("L132").Value = 25650690
$ws.Range("M132").Value = -72989312
$ws.Range("N132").Value = -25655750

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1910.5
$ws.Range("I86").Value = 1924.4788
$ws.Range("K86").Value = 1924.4788
$ws.Range("M86").Value = -801.4788000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 1910.5
$ws.Range("I89").Value = 1924.4788
$ws.Range("K89").Value = 9622.394
$ws.Range("M89").Value = -4006.394

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 770057.4399999999
$ws.Range("I107").Value = 909795.2
$ws.Range("K107").Value = 909795.2
$ws.Range("M107").Value = -907875.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 11292254
$ws.Range("I134").Value = 15724386
$ws.Range("J134").Value = 41454.92
$ws.Range("K134").Value = 47173158
$ws.Range("L134").Value = 124364.76
$ws.Range("M134").Value = -47170623
$ws.Range("N134").Value = -129434.76

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1605723.2
$ws.Range("I31").Value = 2605445
$ws.Range("J31").Value = 6168.8
$ws.Range("K31").Value = 2605445
$ws.Range("L31").Value = 6168.8
$ws.Range("M31").Value = -2605150
$ws.Range("N31").Value = -6758.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1605723.2
$ws.Range("I34").Value = 2605445
$ws.Range("J34").Value = 6168.8
$ws.Range("K34").Value = 2605445
$ws.Range("L34").Value = 6168.8
$ws.Range("M34").Value = -2605243
$ws.Range("N34").Value = -6572.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1340938.5
$ws.Range("I58").Value = 5640.952
$ws.Range("J58").Value = 3497957.8
$ws.Range("K58").Value = 5640.952
$ws.Range("L58").Value = 3497957.8
$ws.Range("M58").Value = -5437.952
$ws.Range("N58").Value = -3498363.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 10550.542
$ws.Range("I99").Value = 7131.769
$ws.Range("J99").Value = 14590.909
$ws.Range("K99").Value = 7131.769
$ws.Range("L99").Value = 14590.909
$ws.Range("M99").Value = -5633.769
$ws.Range("N99").Value = -17586.909

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 565.2692
$ws.Range("I107").Value = 264.81818
$ws.Range("J107").Value = 785.6
$ws.Range("K107").Value = 264.81818
$ws.Range("L107").Value = 785.6
$ws.Range("M107").Value = 1655.18182
$ws.Range("N107").Value = -4625.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 4367.294
$ws.Range("I122").Value = 5107.28
$ws.Range("J122").Value = 2311.7778
$ws.Range("K122").Value = 15321.84
$ws.Range("L122").Value = 6935.3334
$ws.Range("M122").Value = -12871.84
$ws.Range("N122").Value = -11835.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 10550.542
$ws.Range("I126").Value = 7131.769
$ws.Range("J126").Value = 14590.909
$ws.Range("K126").Value = 21395.307
$ws.Range("L126").Value = 43772.727
$ws.Range("M126").Value = -18925.307
$ws.Range("N126").Value = -48712.727

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1484511.1
$ws.Range("I134").Value = 2014.0555
$ws.Range("J134").Value = 4449505.5
$ws.Range("K134").Value = 6042.166499999999
$ws.Range("L134").Value = 13348516.5
$ws.Range("M134").Value = -3507.166499999999
$ws.Range("N134").Value = -13353586.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1340938.5
$ws.Range("I136").Value = 5640.952
$ws.Range("J136").Value = 3497957.8
$ws.Range("K136").Value = 16922.856
$ws.Range("L136").Value = 10493873.4
$ws.Range("M136").Value = -14372.856
$ws.Range("N136").Value = -10498973.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H141").Value = 94213.95
$ws.Range("I141").Value = 29796
$ws.Range("J141").Value = 110318.44
$ws.Range("K141").Value = 29796
$ws.Range("L141").Value = 110318.44
$ws.Range("M141").Value = -24616
$ws.Range("N141").Value = -120678.44

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 6068157
$ws.Range("I121").Value = 388
$ws.Range("J121").Value = 9535453
$ws.Range("K121").Value = 1164
$ws.Range("L121").Value = 28606359
$ws.Range("M121").Value = 146
$ws.Range("N121").Value = -28608979

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 15704.294
$ws.Range("I131").Value = 143322.86
$ws.Range("J131").Value = 1059.541
$ws.Range("K131").Value = 429968.58
$ws.Range("L131").Value = 3178.623
$ws.Range("M131").Value = -424928.58
$ws.Range("N131").Value = -13258.623

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 36.6
$ws.Range("I2").Value = 36.6
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 36.6
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 76.40000000000001
$ws.Range("N2").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 31250146
$ws.Range("I55").Value = 50000056
$ws.Range("J55").Value = 293.33334
$ws.Range("K55").Value = 50000056
$ws.Range("L55").Value = 293.33334
$ws.Range("M55").Value = -49999883
$ws.Range("N55").Value = -639.33334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3246.8
$ws.Range("I61").Value = 2438.8333
$ws.Range("J61").Value = 4458.75
$ws.Range("K61").Value = 2438.8333
$ws.Range("L61").Value = 4458.75
$ws.Range("M61").Value = -2236.8333
$ws.Range("N61").Value = -4862.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 3246.8
$ws.Range("I113").Value = 2438.8333
$ws.Range("J113").Value = 4458.75
$ws.Range("K113").Value = 2438.8333
$ws.Range("L113").Value = 4458.75
$ws.Range("M113").Value = -268.8332999999998
$ws.Range("N113").Value = -8798.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3764089
$ws.Range("I132").Value = 5719387.5
$ws.Range("J132").Value = 3899.923
$ws.Range("K132").Value = 17158162.5
$ws.Range("L132").Value = 11699.769
$ws.Range("M132").Value = -17155632.5
$ws.Range("N132").Value = -16759.769

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 974783.9
$ws.Range("I132").Value = 3106.889
$ws.Range("J132").Value = 3665581.5
$ws.Range("K132").Value = 9320.667000000001
$ws.Range("L132").Value = 10996744.5
$ws.Range("M132").Value = -6790.667000000001
$ws.Range("N132").Value = -11001804.5
